$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4717948717948718
$ws.Range("E2").Value = 0.7076923076923077
$ws.Range("F2").Value = 0.8487179487179487
$ws.Range("G2").Value = 0.7384615384615385
$ws.Range("H2").Value = 0.8683760683760684
$ws.Range("I2").Value = 0.6417910447761194
$ws.Range("K2").Value = 0.7089552238805971
$ws.Range("L2").Value = 0.487012987012987
$ws.Range("N2").Value = 0.6404109589041096
$ws.Range("P2").Value = 0.53125
$ws.Range("Q2").Value = 0.3181818181818182
$ws.Range("R2").Value = 0.3966942148760331
$ws.Range("S2").Value = 0.4615384615384616
$ws.Range("T2").Value = 0.4285714285714285
$ws.Range("V2").Value = 0.4871794871794872
$ws.Range("W2").Value = 0.4418604651162791
$ws.Range("X2").Value = 0.3333333333333333
$ws.Range("Y2").Value = 0.3703703703703703
$ws.Range("AA2").Value = 0.5
